# Auto-generated edit script to apply cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.154.19"
$ws.Range("E2").Value = "  -2.18%  "

# Row 3
$ws.Range("D3").Value = "1.852.48"
$ws.Range("E3").Value = "  -1.00%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.43"
$ws.Range("E5").Value = "  -1.55%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6884"
$ws.Range("E6").Value = "  -5.23%  "

# Row 7
$ws.Range("E7").Value = "  +0.02%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07703"
$ws.Range("E8").Value = "  +7.92%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3038"
$ws.Range("E9").Value = "  -2.99%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.16"
$ws.Range("E10").Value = "  -5.20%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08160"
$ws.Range("E11").Value = "  +0.47%  "

# Row 12
$ws.Range("D12").Value = "1.848.55"
$ws.Range("E12").Value = "  -1.04%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7232"
$ws.Range("E13").Value = "  -2.51%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.198"
$ws.Range("E14").Value = "  -2.76%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.06"
$ws.Range("E15").Value = "  -3.65%  "

# Row 16
$ws.Range("D16").Value = "29.162.21"
$ws.Range("E16").Value = "  -2.15%  "

# Row 17
$ws.Range("B17").Value = "Uniswap"
$ws.Range("C17").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.731"
$ws.Range("E17").Value = "  -4.30%  "

# Row 18
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007793"
$ws.Range("E18").Value = "  +0.00%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.17"
$ws.Range("E19").Value = "  -1.55%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "234.31"
$ws.Range("E20").Value = "  -5.14%  "

# Row 21
$ws.Range("E21").Value = "  +0.07%  "

# Row 22
$ws.Range("D22").Value = "2.101.29"
$ws.Range("E22").Value = "  -1.34%  "

# Row 23
$ws.Range("E23").Value = "  -0.03%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.511"
$ws.Range("E24").Value = "  -2.83%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "161.78"
$ws.Range("E25").Value = "  -1.23%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.956"
$ws.Range("E26").Value = "  -2.83%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1424"
$ws.Range("E27").Value = "  -7.07%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.07"
$ws.Range("E28").Value = "  -2.52%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.960"
$ws.Range("E29").Value = "  -2.40%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.402"
$ws.Range("E30").Value = "  -3.08%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.516"
$ws.Range("E31").Value = "  -0.10%  "

# Row 32
$ws.Range("E32").Value = "  -2.55%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.003"
$ws.Range("E33").Value = "  -4.20%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05188"
$ws.Range("E34").Value = "  -2.51%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.177"
$ws.Range("E35").Value = "  -4.19%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7024"
$ws.Range("E36").Value = "  -4.70%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.020"
$ws.Range("E37").Value = "  +2.21%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.652"
$ws.Range("E38").Value = "  -1.69%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01850"
$ws.Range("E39").Value = "  -4.16%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.678"
$ws.Range("E40").Value = "  -1.99%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9107"
$ws.Range("E41").Value = "  +1.76%  "

# Row 42
$ws.Range("D42").Value = "1.094.25"
$ws.Range("E42").Value = "  +5.22%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.969"
$ws.Range("E43").Value = "  +0.20%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4273"
$ws.Range("E44").Value = "  -4.36%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "70.39"
$ws.Range("E45").Value = "  -1.11%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.001"
$ws.Range("E46").Value = "  -0.01%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.49"
$ws.Range("E47").Value = "  -1.17%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.755"
$ws.Range("E48").Value = "  -3.56%  "

# Row 49
$ws.Range("D49").Value = "1.999.30"
$ws.Range("E49").Value = "  -1.32%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.120"
$ws.Range("E50").Value = "  -4.48%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.917"
$ws.Range("E51").Value = "  -7.36%  "

